$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Drop the first 4 data rows (old Cutoff = 0..3); remaining rows shift up
    # and Excel automatically renumbers the Cutoff/Reaction_number/count
    # columns along with the row shift.
    $ws.Range("A2:C5").EntireRow.Delete()

    # Renumber the "Cutoff" column (A) back to a 0-based sequence now that
    # the top rows are gone.
    for ($i = 0; $i -le 14; $i++) {
        $ws.Cells.Item($i + 2, 1).Value = $i
    }
}
